$wb = $excel.ActiveWorkbook

# --- Neodymium sheet ---
$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C2").Value = [double]"1.624997844765925E-06"
$ws.Range("D2").Value = 0.004327962216169592
$ws.Range("E2").Value = 0.004954444338515564

$ws.Range("B3").Value = [double]"2.183968937109617E-12"
$ws.Range("C3").Value = [double]"7.863771111478995E-05"
$ws.Range("D3").Value = 0.003773862252567578
$ws.Range("E3").Value = 0.004405377924640967

$ws.Range("B4").Value = [double]"3.409239115768776E-14"
$ws.Range("C4").Value = [double]"7.112191862264318E-05"
$ws.Range("D4").Value = 0.003084289105169288
$ws.Range("E4").Value = 0.003890045279137931

$ws.Range("C5").Value = [double]"1.583093649904029E-09"
$ws.Range("D5").Value = 0.0001705074686156489
$ws.Range("E5").Value = 0.0003248552659786882

# --- Copper sheet ---
$ws = $wb.Worksheets.Item("Copper")
$ws.Range("B2").Value = [double]"3.278472098474135E-06"
$ws.Range("C2").Value = 0.0025272965822567
$ws.Range("D2").Value = 0.3647353462301753
$ws.Range("E2").Value = 0.3329926493422712

$ws.Range("B3").Value = [double]"2.22924718813326E-05"
$ws.Range("C3").Value = 0.00911755121305614
$ws.Range("D3").Value = 0.2590053605817614
$ws.Range("E3").Value = 0.2556986289008428

$ws.Range("B4").Value = [double]"6.611256234481376E-05"
$ws.Range("C4").Value = 0.002440576610812275
$ws.Range("D4").Value = 0.2198080158901728
$ws.Range("E4").Value = 0.2575291483701416

$ws.Range("B5").Value = [double]"2.076903987060008E-05"
$ws.Range("C5").Value = 0.005351015551773063
$ws.Range("D5").Value = 0.3156820730433227
$ws.Range("E5").Value = 0.2600971549971024

# --- Raw silicon sheet ---
$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("B2").Value = [double]"4.966311329314386E-07"
$ws.Range("C2").Value = [double]"3.454741537111926E-05"
$ws.Range("D2").Value = 0.009305532546032724
$ws.Range("E2").Value = 0.00862110342949683

$ws.Range("B3").Value = [double]"5.299988190966853E-07"
$ws.Range("C3").Value = 0.0001154307604995803
$ws.Range("D3").Value = 0.004900770646459259
$ws.Range("E3").Value = 0.004790782953538905

$ws.Range("B4").Value = [double]"3.396088080967769E-06"
$ws.Range("C4").Value = [double]"3.239926736874637E-05"
$ws.Range("D4").Value = 0.005064313401933368
$ws.Range("E4").Value = 0.006014344550262347

$ws.Range("B5").Value = [double]"1.82357813169359E-06"
$ws.Range("C5").Value = [double]"4.11437408643171E-05"
$ws.Range("D5").Value = 0.008674230610438876
$ws.Range("E5").Value = 0.007152240119196929
